$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Update the Weight column (D2:D11) with new values
# (order matches the shared-string insertion order of the target workbook)
$ws.Range("D2").Value = "164g"
$ws.Range("D8").Value = "113g"
$ws.Range("D10").Value = "157g"
$ws.Range("D6").Value = "165g"
$ws.Range("D7").Value = "175g"
$ws.Range("D4").Value = "192g"
$ws.Range("D9").Value = "185g"
$ws.Range("D11").Value = "162g"

# Update window / view state: zoom to 235% and move the selection
$win = $excel.ActiveWindow
$win.Zoom = 235
$ws.Range("E14").Select()
